# Season_Attack/89.xlsx edit:
# Add a new date-pair of columns (DV = "07-26_A", DW = "07-26_0") mirroring
# the existing DT ("07-25_A") / DU ("07-25_0") pair, and convert the former
# DU text values into real numbers (matching what DT/DS etc already do).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# 1. Duplicate the DT:DU block (values + formatting) into DV:DW.
#    This gives DV an exact copy of DT (style + value) and DW an exact
#    copy of DU (the original inline-string value), matching the pattern
#    already used by every previous date pair in the sheet.
$srcRange = $ws.Range("DT1:DU$lastRow")
$dstCell = $ws.Range("DV1")
$srcRange.Copy($dstCell)

# 2. Row 1 holds column headers - the new pair gets its own labels
#    instead of being a literal copy of the DT1/DU1 headers.
$ws.Range("DV1").Value = "07-26_A"
$ws.Range("DW1").Value = "07-26_0"

# 3. Convert the old DU column from inline-string numbers to real numeric
#    values for every data row (row 1 is header text, so skip it).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 125)   # column 125 = DU
    if ($cell.Value2 -ne "") {
        $cell.Value = $cell.Value2
    }
}
